$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 2).Value = "Bitcoin"
$ws.Cells.Item(2, 3).Value = "https://coinranking.com/coin/Qwsogvtv82FCd+bitcoin-btc"
$ws.Cells.Item(2, 4).NumberFormat = "@"
$ws.Cells.Item(2, 4).Value = "24.893.62"
$ws.Cells.Item(2, 5).Value = "  +0.49%  "
$ws.Cells.Item(3, 2).Value = "Ethereum"
$ws.Cells.Item(3, 3).Value = "https://coinranking.com/coin/razxDUgYGNAdQ+ethereum-eth"
$ws.Cells.Item(3, 4).NumberFormat = "@"
$ws.Cells.Item(3, 4).Value = "1.661.58"
$ws.Cells.Item(3, 5).Value = "  +0.46%  "
$ws.Cells.Item(4, 2).Value = "TetherUSD"
$ws.Cells.Item(4, 3).Value = "https://coinranking.com/coin/HIVsRcGKkPFtW+tetherusd-usdt"
$ws.Cells.Item(4, 4).NumberFormat = "@"
$ws.Cells.Item(4, 4).Value = "1.003"
$ws.Cells.Item(4, 5).Value = "  -0.43%  "
$ws.Cells.Item(5, 2).Value = "BNB"
$ws.Cells.Item(5, 3).Value = "https://coinranking.com/coin/WcwrkfNI4FUAe+bnb-bnb"
$ws.Cells.Item(5, 4).NumberFormat = "@"
$ws.Cells.Item(5, 4).Value = "325.07"
$ws.Cells.Item(5, 5).Value = "  +5.18%  "
$ws.Cells.Item(6, 2).Value = "USDC"
$ws.Cells.Item(6, 3).Value = "https://coinranking.com/coin/aKzUVe4Hh_CON+usdc-usdc"
$ws.Cells.Item(6, 4).NumberFormat = "@"
$ws.Cells.Item(6, 4).Value = "1.001"
$ws.Cells.Item(6, 5).Value = "  -0.11%  "
$ws.Cells.Item(7, 2).Value = "XRP"
$ws.Cells.Item(7, 4).NumberFormat = "@"
$ws.Cells.Item(7, 4).Value = "0.3632"
$ws.Cells.Item(7, 5).Value = "  -0.32%  "
$ws.Cells.Item(8, 2).Value = "OKB"
$ws.Cells.Item(8, 4).NumberFormat = "@"
$ws.Cells.Item(8, 4).Value = "47.54"
$ws.Cells.Item(8, 5).Value = "  +1.35%  "
$ws.Cells.Item(9, 2).Value = "Cardano"
$ws.Cells.Item(9, 5).Value = "  -1.53%  "
$ws.Cells.Item(10, 2).Value = "Polygon"
$ws.Cells.Item(10, 4).NumberFormat = "@"
$ws.Cells.Item(10, 4).Value = "1.133"
$ws.Cells.Item(10, 5).Value = "  -1.02%  "
$ws.Cells.Item(11, 2).Value = "Dogecoin"
$ws.Cells.Item(11, 4).NumberFormat = "@"
$ws.Cells.Item(11, 4).Value = "0.07081"
$ws.Cells.Item(11, 5).Value = "  -1.52%  "
$ws.Cells.Item(12, 2).Value = "BinanceUSD"
$ws.Cells.Item(12, 5).Value = "  -0.29%  "
$ws.Cells.Item(13, 2).Value = "Polkadot"
$ws.Cells.Item(13, 4).NumberFormat = "@"
$ws.Cells.Item(13, 4).Value = "6.043"
$ws.Cells.Item(13, 5).Value = "  -0.85%  "
$ws.Cells.Item(14, 2).Value = "Solana"
$ws.Cells.Item(14, 4).NumberFormat = "@"
$ws.Cells.Item(14, 4).Value = "19.51"
$ws.Cells.Item(14, 5).Value = "  -2.52%  "
$ws.Cells.Item(15, 2).Value = "WrappedEther"
$ws.Cells.Item(15, 4).NumberFormat = "@"
$ws.Cells.Item(15, 4).Value = "1.656.14"
$ws.Cells.Item(15, 5).Value = "  +0.29%  "
$ws.Cells.Item(16, 2).Value = "Chainlink"
$ws.Cells.Item(16, 4).NumberFormat = "@"
$ws.Cells.Item(16, 4).Value = "6.603"
$ws.Cells.Item(16, 5).Value = "  -1.63%  "
$ws.Cells.Item(17, 2).Value = "ShibaInu"
$ws.Cells.Item(17, 4).NumberFormat = "@"
$ws.Cells.Item(17, 4).Value = "0.00001046"
$ws.Cells.Item(17, 5).Value = "  -2.52%  "
$ws.Cells.Item(18, 2).Value = "TRON"
$ws.Cells.Item(18, 4).NumberFormat = "@"
$ws.Cells.Item(18, 4).Value = "0.06594"
$ws.Cells.Item(18, 5).Value = "  -0.11%  "
$ws.Cells.Item(19, 2).Value = "Dai"
$ws.Cells.Item(19, 4).NumberFormat = "@"
$ws.Cells.Item(19, 4).Value = "1.000"
$ws.Cells.Item(19, 5).Value = "  -0.13%  "
$ws.Cells.Item(20, 2).Value = "Litecoin"
$ws.Cells.Item(20, 4).NumberFormat = "@"
$ws.Cells.Item(20, 4).Value = "79.08"
$ws.Cells.Item(20, 5).Value = "  -1.82%  "
$ws.Cells.Item(21, 2).Value = "Uniswap"
$ws.Cells.Item(21, 4).NumberFormat = "@"
$ws.Cells.Item(21, 4).Value = "5.919"
$ws.Cells.Item(21, 5).Value = "  -2.27%  "
$ws.Cells.Item(22, 2).Value = "Avalanche"
$ws.Cells.Item(22, 4).NumberFormat = "@"
$ws.Cells.Item(22, 4).Value = "15.77"
$ws.Cells.Item(22, 5).Value = "  -4.44%  "
$ws.Cells.Item(23, 2).Value = "Cosmos"
$ws.Cells.Item(23, 4).NumberFormat = "@"
$ws.Cells.Item(23, 4).Value = "12.58"
$ws.Cells.Item(23, 5).Value = "  +2.59%  "
$ws.Cells.Item(24, 2).Value = "WrappedBTC"
$ws.Cells.Item(24, 4).NumberFormat = "@"
$ws.Cells.Item(24, 4).Value = "24.856.77"
$ws.Cells.Item(24, 5).Value = "  +0.40%  "
$ws.Cells.Item(25, 2).Value = "Toncoin"
$ws.Cells.Item(25, 4).NumberFormat = "@"
$ws.Cells.Item(25, 4).Value = "2.449"
$ws.Cells.Item(25, 5).Value = "  +1.30%  "
$ws.Cells.Item(26, 2).Value = "LidoDAOToken"
$ws.Cells.Item(26, 4).NumberFormat = "@"
$ws.Cells.Item(26, 4).Value = "2.430"
$ws.Cells.Item(26, 5).Value = "  -5.69%  "
$ws.Cells.Item(27, 2).Value = "Monero"
$ws.Cells.Item(27, 4).NumberFormat = "@"
$ws.Cells.Item(27, 4).Value = "148.48"
$ws.Cells.Item(27, 5).Value = "  -1.15%  "
$ws.Cells.Item(28, 2).Value = "EthereumClassic"
$ws.Cells.Item(28, 5).Value = "  -4.17%  "
$ws.Cells.Item(29, 2).Value = "WrappedliquidstakedEther2.0"
$ws.Cells.Item(29, 4).NumberFormat = "@"
$ws.Cells.Item(29, 4).Value = "1.835.48"
$ws.Cells.Item(29, 5).Value = "  -0.20%  "
$ws.Cells.Item(30, 2).Value = "BitcoinCash"
$ws.Cells.Item(30, 4).NumberFormat = "@"
$ws.Cells.Item(30, 4).Value = "125.33"
$ws.Cells.Item(30, 5).Value = "  -2.31%  "
$ws.Cells.Item(31, 2).Value = "ImmutableX"
$ws.Cells.Item(31, 4).NumberFormat = "@"
$ws.Cells.Item(31, 4).Value = "1.186"
$ws.Cells.Item(31, 5).Value = "  -1.94%  "
$ws.Cells.Item(32, 2).Value = "HuobiToken"
$ws.Cells.Item(32, 4).NumberFormat = "@"
$ws.Cells.Item(32, 4).Value = "4.085"
$ws.Cells.Item(32, 5).Value = "  -1.12%  "
$ws.Cells.Item(33, 2).Value = "Filecoin"
$ws.Cells.Item(33, 4).NumberFormat = "@"
$ws.Cells.Item(33, 4).Value = "5.735"
$ws.Cells.Item(33, 5).Value = "  -7.74%  "
$ws.Cells.Item(34, 2).Value = "Stellar"
$ws.Cells.Item(34, 4).NumberFormat = "@"
$ws.Cells.Item(34, 4).Value = "0.08473"
$ws.Cells.Item(34, 5).Value = "  -0.34%  "
$ws.Cells.Item(35, 2).Value = "WEMIXTOKEN"
$ws.Cells.Item(35, 4).NumberFormat = "@"
$ws.Cells.Item(35, 4).Value = "1.645"
$ws.Cells.Item(35, 5).Value = "  -4.47%  "
$ws.Cells.Item(36, 2).Value = "Aptos"
$ws.Cells.Item(36, 4).NumberFormat = "@"
$ws.Cells.Item(36, 4).Value = "12.19"
$ws.Cells.Item(36, 5).Value = "  -5.89%  "
$ws.Cells.Item(37, 2).Value = "TrustWalletToken"
$ws.Cells.Item(37, 4).NumberFormat = "@"
$ws.Cells.Item(37, 4).Value = "1.284"
$ws.Cells.Item(37, 5).Value = "  +4.23%  "
$ws.Cells.Item(38, 2).Value = "InternetComputer(DFINITY)"
$ws.Cells.Item(38, 4).NumberFormat = "@"
$ws.Cells.Item(38, 4).Value = "5.166"
$ws.Cells.Item(38, 5).Value = "  -2.29%  "
$ws.Cells.Item(39, 2).Value = "VeChain"
$ws.Cells.Item(39, 4).NumberFormat = "@"
$ws.Cells.Item(39, 4).Value = "0.02272"
$ws.Cells.Item(39, 5).Value = "  -1.40%  "
$ws.Cells.Item(40, 2).Value = "Hedera"
$ws.Cells.Item(40, 4).NumberFormat = "@"
$ws.Cells.Item(40, 4).Value = "0.06102"
$ws.Cells.Item(40, 5).Value = "  -2.82%  "
$ws.Cells.Item(41, 2).Value = "FraxShare"
$ws.Cells.Item(41, 4).NumberFormat = "@"
$ws.Cells.Item(41, 4).Value = "8.310"
$ws.Cells.Item(41, 5).Value = "  -1.44%  "
$ws.Cells.Item(42, 2).Value = "Algorand"
$ws.Cells.Item(42, 4).NumberFormat = "@"
$ws.Cells.Item(42, 4).Value = "0.2074"
$ws.Cells.Item(42, 5).Value = "  -2.07%  "
$ws.Cells.Item(43, 2).Value = "Frax"
$ws.Cells.Item(43, 5).Value = "  -0.01%  "
$ws.Cells.Item(44, 2).Value = "TheSandbox"
$ws.Cells.Item(44, 4).NumberFormat = "@"
$ws.Cells.Item(44, 4).Value = "0.5925"
$ws.Cells.Item(44, 5).Value = "  -3.54%  "
$ws.Cells.Item(45, 2).Value = "PancakeSwap"
$ws.Cells.Item(45, 3).Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Cells.Item(45, 4).NumberFormat = "@"
$ws.Cells.Item(45, 4).Value = "3.830"
$ws.Cells.Item(45, 5).Value = "  +1.70%  "
$ws.Cells.Item(46, 2).Value = "EnergySwap"
$ws.Cells.Item(46, 3).Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Cells.Item(46, 4).NumberFormat = "@"
$ws.Cells.Item(46, 4).Value = "13.23"
$ws.Cells.Item(46, 5).Value = "  -0.24%  "
$ws.Cells.Item(47, 2).Value = "Decentraland"
$ws.Cells.Item(47, 4).NumberFormat = "@"
$ws.Cells.Item(47, 4).Value = "0.5623"
$ws.Cells.Item(47, 5).Value = "  -3.77%  "
$ws.Cells.Item(48, 2).Value = "Quant"
$ws.Cells.Item(48, 4).NumberFormat = "@"
$ws.Cells.Item(48, 4).Value = "125.02"
$ws.Cells.Item(48, 5).Value = "  +0.70%  "
$ws.Cells.Item(49, 2).Value = "NEARProtocol"
$ws.Cells.Item(49, 4).NumberFormat = "@"
$ws.Cells.Item(49, 4).Value = "1.946"
$ws.Cells.Item(49, 5).Value = "  -3.01%  "
$ws.Cells.Item(50, 2).Value = "Cronos"
$ws.Cells.Item(50, 4).NumberFormat = "@"
$ws.Cells.Item(50, 4).Value = "0.06986"
$ws.Cells.Item(50, 5).Value = "  -1.30%  "
$ws.Cells.Item(51, 2).Value = "Tezos"
$ws.Cells.Item(51, 4).NumberFormat = "@"
$ws.Cells.Item(51, 4).Value = "1.188"
$ws.Cells.Item(51, 5).Value = "  +0.58%  "
